$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Materials" section header, styled like the existing "Hardware" header (A12)
$ws.Range("A17").Value = "Materials"
$ws.Range("A17").Font.Color = $ws.Range("A12").Font.Color
$ws.Range("A17").Font.Bold = $ws.Range("A12").Font.Bold

# New material line item
$ws.Range("A18").Value = "1/16"" (1.6mm) acrylic for laser-cut buffer layer"
$ws.Range("B18").Value = "1 buffer per link"

# Match the saved selection/active cell from the source workbook
$ws.Range("C17").Select()
